# Generate Report for Handback
# Updates the timestamps (and one "ht" -> "mt" status value) recorded in the
# handback-status workbook to reflect a newly re-run report generation.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet -------------------------------------------------
# Latest HO Xliff Generate Date for the 2276b118-... file, shown on both
# the row that references it directly (row 2) and the row that shares the
# same generation timestamp (row 5).
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 00:27:48"
$wsOverview.Range("G5").Value = "2016-09-07 00:27:48"

# --- "zh-cn" sheet ------------------------------------------------------
$wsZh = $wb.Sheets.Item("zh-cn")
# Priority changed from "ht" to "mt"
$wsZh.Range("E2").Value = "mt"
$wsZh.Range("E5").Value = "mt"
# Correspond Handoff Datetime
$wsZh.Range("H2").Value = "2016-09-07 00:27:43"
$wsZh.Range("H5").Value = "2016-09-07 00:27:43"
# Correspond Handback DateTime
$wsZh.Range("K2").Value = "2016-09-07 00:28:09"
$wsZh.Range("K5").Value = "2016-09-07 00:28:09"

# --- "de-de" sheet ------------------------------------------------------
$wsDe = $wb.Sheets.Item("de-de")
# Priority changed from "ht" to "mt"
$wsDe.Range("E2").Value = "mt"
$wsDe.Range("E5").Value = "mt"
# Correspond Handoff Datetime (shares the same value as Overview's
# "Latest HO Xliff Generate Date" for this file)
$wsDe.Range("H2").Value = "2016-09-07 00:27:48"
$wsDe.Range("H5").Value = "2016-09-07 00:27:48"
# Correspond Handback DateTime
$wsDe.Range("K2").Value = "2016-09-07 00:28:17"
$wsDe.Range("K5").Value = "2016-09-07 00:28:17"
